# Refresh the cryptocurrency price / 1h-volume-change table with the
# latest values from the coinranking.com scrape (GitHub Actions job).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.854.20"
$ws.Range("E2").Value = "  +1.30%  "

$ws.Range("D3").Value = "2.621.58"
$ws.Range("E3").Value = "  +1.10%  "

$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").Value = "'603.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.51%  "

$ws.Range("D6").Value = "'155.05"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.14%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").Value = "'0.548"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.72%  "

$ws.Range("D9").Value = "2.621.35"
$ws.Range("E9").Value = "  +1.11%  "

$ws.Range("D10").Value = "'0.127"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +9.31%  "

$ws.Range("E11").Value = "  +1.00%  "

$ws.Range("E12").Value = "  +0.72%  "

$ws.Range("E13").Value = "  -1.46%  "

$ws.Range("D14").Value = "'27.70"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.88%  "

$ws.Range("E15").Value = "  +2.83%  "

$ws.Range("D16").Value = "3.101.03"
$ws.Range("E16").Value = "  +1.50%  "

$ws.Range("D17").Value = "67.808.48"
$ws.Range("E17").Value = "  +1.11%  "

$ws.Range("D18").Value = "2.626.63"

$ws.Range("D19").Value = "'11.20"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.20%  "

$ws.Range("D20").Value = "'366.76"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.86%  "

$ws.Range("D21").Value = "'7.65"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.41%  "

$ws.Range("E22").Value = "  -0.72%  "

$ws.Range("E23").Value = "  -2.51%  "

$ws.Range("E24").Value = "  +0.11%  "

$ws.Range("D25").Value = "'9.86"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -7.33%  "

$ws.Range("D26").Value = "'66.26"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.49%  "

$ws.Range("E28").Value = "  -0.48%  "

$ws.Range("D29").Value = "'577.49"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.77%  "

$ws.Range("D30").Value = "'1.02"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.96%  "

$ws.Range("E31").Value = "  -3.22%  "

$ws.Range("E32").Value = "  -2.45%  "

$ws.Range("D33").Value = "'1.87"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.68%  "

$ws.Range("E34").Value = "  -1.74%  "

$ws.Range("E35").Value = "  +0.04%  "

$ws.Range("E36").Value = "  -3.88%  "

$ws.Range("E37").Value = "  -2.36%  "

$ws.Range("D38").Value = "'158.22"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.39%  "

$ws.Range("D39").Value = "'19.37"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.26%  "

$ws.Range("E40").Value = "  -0.42%  "

$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D41").Value = "'5.35"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.95%  "

$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "'1.85"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.67%  "

$ws.Range("D43").Value = "'2.58"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.92%  "

$ws.Range("D44").Value = "'41.22"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.89%  "

$ws.Range("E46").Value = "  -0.07%  "

$ws.Range("D47").Value = "'157.14"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.34%  "

$ws.Range("D48").Value = "0.0₆0286"
$ws.Range("E48").Value = "  -7.84%  "

$ws.Range("E49").Value = "  -0.58%  "

$ws.Range("E50").Value = "  -2.19%  "

$ws.Range("D51").Value = "'0.624"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.45%  "
